{"js": "// Update the date heading (first paragraph of the document body). Using\n// insertText(..., replace) on the paragraph keeps its existing run\n// formatting (font/size) intact.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.insertText(\"2023-08-10 Thursday\", Word.InsertLocation.replace);\n\n// Update every arithmetic-problem cell in the 20x5 practice table, in\n// document order (row-major). Assigning Cell.value (rather than\n// replacing the cell body's text) preserves each cell's existing\n// paragraph/run formatting (alignment, font, size).\nconst newValues = [\n  [\"74-26=\", \"47+25=\", \"54-46=\", \"43-25=\", \"19+13=\"],\n  [\"84-68=\", \"17+35=\", \"39+18=\", \"84-27=\", \"49+45=\"],\n  [\"73-15=\", \"18+14=\", \"91-29=\", \"28-9=\", \"71-65=\"],\n  [\"83-39=\", \"44-29=\", \"82-79=\", \"16+15=\", \"47+15=\"],\n  [\"26+19=\", \"25+38=\", \"57-18=\", \"75-67=\", \"26+15=\"],\n  [\"12+49=\", \"61-55=\", \"6+37=\", \"64-19=\", \"31-8=\"],\n  [\"17+44=\", \"66+18=\", \"43-9=\", \"27-19=\", \"33+59=\"],\n  [\"42+9=\", \"16+57=\", \"9+5=\", \"19+65=\", \"19+34=\"],\n  [\"69+9=\", \"36+15=\", \"22-9=\", \"91-6=\", \"28+24=\"],\n  [\"35-7=\", \"34-8=\", \"66+6=\", \"29+63=\", \"35-7=\"],\n  [\"40-3=\", \"95-9=\", \"38+46=\", \"19+42=\", \"67+8=\"],\n  [\"35-26=\", \"91-24=\", \"19+49=\", \"65-39=\", \"9+55=\"],\n  [\"54+39=\", \"83-49=\", \"41-38=\", \"93-47=\", \"71-43=\"],\n  [\"82-14=\", \"92-18=\", \"27+65=\", \"56-17=\", \"40-7=\"],\n  [\"84-58=\", \"75-47=\", \"72-66=\", \"27+68=\", \"80-23=\"],\n  [\"33-8=\", \"58+3=\", \"18+3=\", \"91-37=\", \"28+39=\"],\n  [\"39+46=\", \"42-5=\", \"5+59=\", \"68+25=\", \"60-5=\"],\n  [\"58+25=\", \"17+69=\", \"50-38=\", \"48+19=\", \"72-37=\"],\n  [\"52-38=\", \"33+9=\", \"82-64=\", \"64+17=\", \"78-29=\"],\n  [\"15-6=\", \"33+8=\", \"60-57=\", \"60-1=\", \"92-55=\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst columnCount = table.values[0].length;\n\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < columnCount; c++) {\n    const cell = table.getCellOrNullObject(r, c);\n    cell.value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading (first paragraph of the document body),\n# preserving its existing run formatting.\n$d = $word.ActiveDocument\n$d.Paragraphs.Item(1).Range.Text = '2023-08-10 Thursday'\n\n# Update every arithmetic-problem cell in the 20x5 practice table, in\n# document order (row-major). Re-fetching Cell() after each write keeps\n# the COM Range reference valid and preserves each cell's run formatting\n# because only the text run content is replaced.\n$newValues = @(\n    @('74-26=', '47+25=', '54-46=', '43-25=', '19+13='),\n    @('84-68=', '17+35=', '39+18=', '84-27=', '49+45='),\n    @('73-15=', '18+14=', '91-29=', '28-9=', '71-65='),\n    @('83-39=', '44-29=', '82-79=', '16+15=', '47+15='),\n    @('26+19=', '25+38=', '57-18=', '75-67=', '26+15='),\n    @('12+49=', '61-55=', '6+37=', '64-19=', '31-8='),\n    @('17+44=', '66+18=', '43-9=', '27-19=', '33+59='),\n    @('42+9=', '16+57=', '9+5=', '19+65=', '19+34='),\n    @('69+9=', '36+15=', '22-9=', '91-6=', '28+24='),\n    @('35-7=', '34-8=', '66+6=', '29+63=', '35-7='),\n    @('40-3=', '95-9=', '38+46=', '19+42=', '67+8='),\n    @('35-26=', '91-24=', '19+49=', '65-39=', '9+55='),\n    @('54+39=', '83-49=', '41-38=', '93-47=', '71-43='),\n    @('82-14=', '92-18=', '27+65=', '56-17=', '40-7='),\n    @('84-58=', '75-47=', '72-66=', '27+68=', '80-23='),\n    @('33-8=', '58+3=', '18+3=', '91-37=', '28+39='),\n    @('39+46=', '42-5=', '5+59=', '68+25=', '60-5='),\n    @('58+25=', '17+69=', '50-38=', '48+19=', '72-37='),\n    @('52-38=', '33+9=', '82-64=', '64+17=', '78-29='),\n    @('15-6=', '33+8=', '60-57=', '60-1=', '92-55=')\n)\n\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $rowValues.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
